# Generate Report for Handoff
# Updates the localization-status workbook: the 89c4ac4b... and
# a2d1323e... source files have moved from "Handed back: in sync with
# en-US" to "Ready for handoff" with a fresh handoff timestamp, and the
# zh-cn / de-de detail sheets now carry an "Error Detail" message noting
# that the handback file version is stale.

$wb = $excel.ActiveWorkbook

$msg89 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0859b048312b51ed9af303bb9af15afdbc32ab9c/e2e/89c4ac4b-6781-4913-a4b4-f4f5eb384af7.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a4e98e28330e86694e3ec6f28a557a30b8af740b/e2e/89c4ac4b-6781-4913-a4b4-f4f5eb384af7.md."
$msgA2 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0859b048312b51ed9af303bb9af15afdbc32ab9c/e2e/a2d1323e-5cca-409e-bf91-de84127ba9b0.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a4e98e28330e86694e3ec6f28a557a30b8af740b/e2e/a2d1323e-5cca-409e-bf91-de84127ba9b0.md."

# ---- Overview sheet (rows 4 & 5 = 89c4ac4b*, a2d1323e*) ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Cells.Item(4, 5).Value = "Ready for handoff"
$ov.Cells.Item(4, 6).Value = "Ready for handoff"
$ov.Cells.Item(4, 7).Value = "2016-08-26 12:26:16"
$ov.Cells.Item(5, 5).Value = "Ready for handoff"
$ov.Cells.Item(5, 6).Value = "Ready for handoff"
$ov.Cells.Item(5, 7).Value = "2016-08-26 12:26:16"

# ---- zh-cn detail sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Cells.Item(4, 3).Value = "Ready for handoff"
$zh.Cells.Item(4, 8).Value = "2016-08-26 12:26:16"
$zh.Cells.Item(4, 16).Value = $msg89
$zh.Cells.Item(5, 3).Value = "Ready for handoff"
$zh.Cells.Item(5, 8).Value = "2016-08-26 12:26:10"
$zh.Cells.Item(5, 16).Value = $msgA2
$zh.Columns.Item(16).ColumnWidth = 39.16666666666667

# ---- de-de detail sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Cells.Item(4, 3).Value = "Ready for handoff"
$de.Cells.Item(4, 8).Value = "2016-08-26 12:26:16"
$de.Cells.Item(4, 16).Value = $msg89
$de.Cells.Item(5, 3).Value = "Ready for handoff"
$de.Cells.Item(5, 8).Value = "2016-08-26 12:26:16"
$de.Cells.Item(5, 16).Value = $msgA2
$de.Columns.Item(16).ColumnWidth = 39.16666666666667
